# "Time Log.xlsx" - log three more days of work (2014-09-26..28) in the
# Sheet1 time table, then leave the view scrolled/selected where the user
# ended up.
#
# Values are written as the raw serial numbers Excel itself stores for
# dates/times (days since 1899-12-30 resp. fraction-of-day) instead of
# .NET DateTime objects, so the existing column number formats (custom
# time format, "N mins", etc.) are preserved instead of Excel guessing a
# brand-new date/time format for the newly-populated cells.
#
# D (Interruption) is written before A/B/C on each row purely so that the
# dependent "Delta" formula in column E recalculates against the final
# row values rather than against a transient D=0 state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 64: 2014-09-26, 16:45 -> 17:32, 10 min interruption, Coding
$ws1.Range("D64").Value = 10
$ws1.Range("A64").Value = 41908
$ws1.Range("B64").Value = 0.69791666666666663
$ws1.Range("C64").Value = 0.73055555555555562
$ws1.Range("F64").Value = "Coding"

# Row 65: 2014-09-27, 15:16 -> 16:17, 5 min interruption, Coding
$ws1.Range("D65").Value = 5
$ws1.Range("A65").Value = 41909
$ws1.Range("B65").Value = 0.63611111111111118
$ws1.Range("C65").Value = 0.67847222222222225
$ws1.Range("F65").Value = "Coding"

# Row 66: 2014-09-28, 22:43 -> 23:59, 15 min interruption, Coding
$ws1.Range("D66").Value = 15
$ws1.Range("A66").Value = 41910
$ws1.Range("B66").Value = 0.94652777777777775
$ws1.Range("C66").Value = 0.99930555555555556
$ws1.Range("F66").Value = "Coding"

# Match where the user scrolled to / what they had selected afterwards.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C67").Select()

$wb.Save()
